$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $ws.Range("B24").Formula = "=SUM(`$D24:`$F24)"
    $ws.Range("B25").Formula = "=SUM(`$D25:`$F25)"
}
